$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73; this shifts existing rows 73-119 down to 74-120,
# carrying along their values and formatting (matching the diff's row renumbering).
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record's data.
$ws.Range("A73").Value = 6
$ws.Range("B73").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C73").Value = "Metropolitana"
$ws.Range("D73").Value = 45205
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 100112035
$ws.Range("G73").Value = "Bruselas (repollito)"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 600
$ws.Range("K73").Value = 16000
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = 17067
$ws.Range("N73").Value = "`$/malla 15 kilos"
$ws.Range("O73").Value = "Provincia de Quillota"
$ws.Range("P73").Value = 1138
$ws.Range("Q73").Value = 15
$ws.Range("R73").Value = "Hortaliza"
